# This script applies the updated NATMI ligand-receptor statistics for the
# Efna5 (ligand) -> Epha1 (receptor) pair after the underlying TPM expression
# matrix was refreshed ("update scripts wuth new tpm").
#
# The workbook has one data sheet with columns:
#   A Sending cluster            B Ligand symbol              C Receptor symbol
#   D Target cluster             E Ligand-expressing cells    F Ligand detection rate
#   G Ligand average expr value  H Ligand total expr value
#   I Ligand avg-expr specificity  J Ligand total-expr specificity
#   K Receptor-expressing cells  L Receptor detection rate
#   M Receptor average expr value  N Receptor total expr value
#   O Receptor avg-expr specificity  P Receptor total-expr specificity
#   Q Edge average expr weight   R Edge total expr weight
#   S Edge average expr specificity  T Edge total expr specificity
#
# With the new TPM matrix, the ligand (Efna5) stats change for the "ECs"
# sending cluster, and the receptor (Epha1) stats change for the "ECs",
# "MuSCs" and "Resolving-Mac" target clusters. Every specificity/edge column
# that is derived from those numbers is updated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; Col=5; Value=3}  # E2
    @{Row=2; Col=6; Value=1}  # F2
    @{Row=2; Col=7; Value=0.2708946666666667}  # G2
    @{Row=2; Col=8; Value=0.812684}  # H2
    @{Row=2; Col=9; Value=0.1616296696421007}  # I2
    @{Row=2; Col=10; Value=0.1616296696421007}  # J2
    @{Row=2; Col=13; Value=5.740110333333334}  # M2
    @{Row=2; Col=14; Value=17.220331}  # N2
    @{Row=2; Col=15; Value=0.2861925343043439}  # O2
    @{Row=2; Col=16; Value=0.2861925343043439}  # P2
    @{Row=2; Col=17; Value=1.554965275378222}  # Q2
    @{Row=2; Col=18; Value=13.994687478404}  # R2
    @{Row=2; Col=19; Value=0.04625720477364667}  # S2
    @{Row=2; Col=20; Value=0.04625720477364666}  # T2
    @{Row=3; Col=5; Value=3}  # E3
    @{Row=3; Col=6; Value=1}  # F3
    @{Row=3; Col=7; Value=0.2708946666666667}  # G3
    @{Row=3; Col=8; Value=0.812684}  # H3
    @{Row=3; Col=9; Value=0.1616296696421007}  # I3
    @{Row=3; Col=10; Value=0.1616296696421007}  # J3
    @{Row=3; Col=15; Value=0.2917347240316885}  # O3
    @{Row=3; Col=16; Value=0.2917347240316885}  # P3
    @{Row=3; Col=17; Value=1.585077565332}  # Q3
    @{Row=3; Col=18; Value=14.265698087988}  # R3
    @{Row=3; Col=19; Value=0.04715298706837122}  # S3
    @{Row=3; Col=20; Value=0.04715298706837122}  # T3
    @{Row=4; Col=5; Value=3}  # E4
    @{Row=4; Col=6; Value=1}  # F4
    @{Row=4; Col=7; Value=0.2708946666666667}  # G4
    @{Row=4; Col=8; Value=0.812684}  # H4
    @{Row=4; Col=9; Value=0.1616296696421007}  # I4
    @{Row=4; Col=10; Value=0.1616296696421007}  # J4
    @{Row=4; Col=13; Value=6.759986}  # M4
    @{Row=4; Col=14; Value=20.279958}  # N4
    @{Row=4; Col=15; Value=0.3370418707750538}  # O4
    @{Row=4; Col=16; Value=0.3370418707750538}  # P4
    @{Row=4; Col=17; Value=1.831244154141334}  # Q4
    @{Row=4; Col=18; Value=16.481197387272}  # R4
    @{Row=4; Col=19; Value=0.05447596622892753}  # S4
    @{Row=4; Col=20; Value=0.05447596622892753}  # T4
    @{Row=5; Col=5; Value=3}  # E5
    @{Row=5; Col=6; Value=1}  # F5
    @{Row=5; Col=7; Value=0.2708946666666667}  # G5
    @{Row=5; Col=8; Value=0.812684}  # H5
    @{Row=5; Col=9; Value=0.1616296696421007}  # I5
    @{Row=5; Col=10; Value=0.1616296696421007}  # J5
    @{Row=5; Col=13; Value=1.705448333333333}  # M5
    @{Row=5; Col=14; Value=5.116345}  # N5
    @{Row=5; Col=15; Value=0.0850308708889137}  # O5
    @{Row=5; Col=16; Value=0.0850308708889137}  # P5
    @{Row=5; Col=17; Value=0.4619968577755555}  # Q5
    @{Row=5; Col=18; Value=4.15797171998}  # R5
    @{Row=5; Col=19; Value=0.01374351157115524}  # S5
    @{Row=5; Col=20; Value=0.01374351157115524}  # T5
    @{Row=6; Col=9; Value=0.6313295261673385}  # I6
    @{Row=6; Col=10; Value=0.6313295261673384}  # J6
    @{Row=6; Col=13; Value=5.740110333333334}  # M6
    @{Row=6; Col=14; Value=17.220331}  # N6
    @{Row=6; Col=15; Value=0.2861925343043439}  # O6
    @{Row=6; Col=16; Value=0.2861925343043439}  # P6
    @{Row=6; Col=17; Value=6.073733199387113}  # Q6
    @{Row=6; Col=18; Value=54.66359879448401}  # R6
    @{Row=6; Col=19; Value=0.1806817970749912}  # S6
    @{Row=6; Col=20; Value=0.1806817970749912}  # T6
    @{Row=7; Col=9; Value=0.6313295261673385}  # I7
    @{Row=7; Col=10; Value=0.6313295261673384}  # J7
    @{Row=7; Col=15; Value=0.2917347240316885}  # O7
    @{Row=7; Col=16; Value=0.2917347240316885}  # P7
    @{Row=7; Col=17; Value=6.191352555972}  # Q7
    @{Row=7; Col=18; Value=55.722173003748}  # R7
    @{Row=7; Col=19; Value=0.1841807450894852}  # S7
    @{Row=7; Col=20; Value=0.1841807450894851}  # T7
    @{Row=8; Col=9; Value=0.6313295261673385}  # I8
    @{Row=8; Col=10; Value=0.6313295261673384}  # J8
    @{Row=8; Col=13; Value=6.759986}  # M8
    @{Row=8; Col=14; Value=20.279958}  # N8
    @{Row=8; Col=15; Value=0.3370418707750538}  # O8
    @{Row=8; Col=16; Value=0.3370418707750538}  # P8
    @{Row=8; Col=17; Value=7.152885399634668}  # Q8
    @{Row=8; Col=18; Value=64.375968596712}  # R8
    @{Row=8; Col=19; Value=0.2127844845749681}  # S8
    @{Row=8; Col=20; Value=0.212784484574968}  # T8
    @{Row=9; Col=9; Value=0.6313295261673385}  # I9
    @{Row=9; Col=10; Value=0.6313295261673384}  # J9
    @{Row=9; Col=13; Value=1.705448333333333}  # M9
    @{Row=9; Col=14; Value=5.116345}  # N9
    @{Row=9; Col=15; Value=0.0850308708889137}  # O9
    @{Row=9; Col=16; Value=0.0850308708889137}  # P9
    @{Row=9; Col=17; Value=1.804571264397778}  # Q9
    @{Row=9; Col=18; Value=16.24114137958}  # R9
    @{Row=9; Col=19; Value=0.05368249942789403}  # S9
    @{Row=9; Col=20; Value=0.05368249942789401}  # T9
    @{Row=10; Col=9; Value=0.2070408041905609}  # I10
    @{Row=10; Col=10; Value=0.2070408041905609}  # J10
    @{Row=10; Col=13; Value=5.740110333333334}  # M10
    @{Row=10; Col=14; Value=17.220331}  # N10
    @{Row=10; Col=15; Value=0.2861925343043439}  # O10
    @{Row=10; Col=16; Value=0.2861925343043439}  # P10
    @{Row=10; Col=17; Value=1.991845072848222}  # Q10
    @{Row=10; Col=18; Value=17.926605655634}  # R10
    @{Row=10; Col=19; Value=0.05925353245570605}  # S10
    @{Row=10; Col=20; Value=0.05925353245570605}  # T10
    @{Row=11; Col=9; Value=0.2070408041905609}  # I11
    @{Row=11; Col=10; Value=0.2070408041905609}  # J11
    @{Row=11; Col=15; Value=0.2917347240316885}  # O11
    @{Row=11; Col=16; Value=0.2917347240316885}  # P11
    @{Row=11; Col=19; Value=0.06040099187383215}  # S11
    @{Row=11; Col=20; Value=0.06040099187383215}  # T11
    @{Row=12; Col=9; Value=0.2070408041905609}  # I12
    @{Row=12; Col=10; Value=0.2070408041905609}  # J12
    @{Row=12; Col=13; Value=6.759986}  # M12
    @{Row=12; Col=14; Value=20.279958}  # N12
    @{Row=12; Col=15; Value=0.3370418707750538}  # O12
    @{Row=12; Col=16; Value=0.3370418707750538}  # P12
    @{Row=12; Col=17; Value=2.345746688601334}  # Q12
    @{Row=12; Col=18; Value=21.111720197412}  # R12
    @{Row=12; Col=19; Value=0.06978141997115825}  # S12
    @{Row=12; Col=20; Value=0.06978141997115825}  # T12
    @{Row=13; Col=9; Value=0.2070408041905609}  # I13
    @{Row=13; Col=10; Value=0.2070408041905609}  # J13
    @{Row=13; Col=13; Value=1.705448333333333}  # M13
    @{Row=13; Col=14; Value=5.116345}  # N13
    @{Row=13; Col=15; Value=0.0850308708889137}  # O13
    @{Row=13; Col=16; Value=0.0850308708889137}  # P13
    @{Row=13; Col=17; Value=0.5917985304255555}  # Q13
    @{Row=13; Col=18; Value=5.326186773830001}  # R13
    @{Row=13; Col=19; Value=0.01760485988986445}  # S13
    @{Row=13; Col=20; Value=0.01760485988986445}  # T13
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}
